# feat: Recipe DB, Excel 추가
# Adds 6 new "Combination" recipe items (rows 14-19) to the Items table,
# resizes the table/autofilter/_FilterDatabase range to A1:G19,
# widens column D, and matches the pasted-in mixed fonts/colors that come
# along with the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# ---- New data rows -------------------------------------------------
$ws.Range("A14").Value = 2001001
$ws.Range("B14").Value = "에메랄드 에센스"
$ws.Range("C14").Value = "에메랄드 빛의 에센스. 향긋한 허브향이 나는 꾸덕한 느낌의 향신료이다."
$ws.Range("D14").Value = "Combination"
$ws.Range("E14").Value = 20
$ws.Range("F14").Value = 100

$ws.Range("A15").Value = 2001002
$ws.Range("B15").Value = "루비 물약"
$ws.Range("C15").Value = "루비 색의 물약. 매콤하면서 강렬한 향이 나는 물약으로 남자들이 주로 많이 찾는다."
$ws.Range("D15").Value = "Combination"
$ws.Range("E15").Value = 20
$ws.Range("F15").Value = 200

$ws.Range("A16").Value = 2001003
$ws.Range("B16").Value = "별가루 캔디"
$ws.Range("C16").Value = "밤하늘의 별을 따온 듯한 모습의 사탕. 입 안에서 달달한 폭죽이 터지면서 먹은 이들의 행복을 찾아준다."
$ws.Range("D16").Value = "Combination"
$ws.Range("E16").Value = 20
$ws.Range("F16").Value = 300

$ws.Range("A17").Value = 2001004
$ws.Range("B17").Value = "핑크 프로틴 쉐이크"
$ws.Range("C17").Value = "헬스인들을 위한 분홍빛의 프로틴 쉐이크. 매일 아침 챙겨 먹으면 건강해 질 것 같다."
$ws.Range("D17").Value = "Combination"
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 400

$ws.Range("A18").Value = 2001005
$ws.Range("B18").Value = "오스틴-코기 포도잼"
$ws.Range("C18").Value = "젤리같은 느낌의 포도잼. 품격있는 오스틴 포도의 맛과 중독성 있는 향이 나는 맛있는 잼이다."
$ws.Range("D18").Value = "Combination"
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = 500

$ws.Range("A19").Value = 2001006
$ws.Range("B19").Value = "민트 스톡"
$ws.Range("C19").Value = "민트의 향이 나는 향신료. 풍부한 감칠맛과 향긋한 민트향으로 주부들에게 인기가 많다."
$ws.Range("D19").Value = "Combination"
$ws.Range("E19").Value = 20
$ws.Range("F19").Value = 600

# ---- Re-apply the wrapping "description" look on the new C column ----
# (new rows don't inherit the C-column wrap/dark-gray formatting automatically)
$descRange = $ws.Range("C14:C19")
$descRange.Font.Color = 1907741
$descRange.WrapText = $true

# ---- Mixed fonts on column B (pasted-in text look) ------------------
# B15 ("루비 물약") carries a distinct Arial run with dark-gray text,
# B16-B18 keep the workbook's "맑은 고딕" family but in the same
# dark-gray color used elsewhere in the sheet.
$ws.Range("B15").Font.Color = 1907741
$ws.Range("B15").Font.Family = 2
$ws.Range("B15").Font.Name = "Arial"

$ws.Range("B16:B18").Font.Color = 1907741

# ---- Table / AutoFilter resize --------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G19"))

# ---- _FilterDatabase defined name ------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Items!_FilterDatabase") {
        $n.RefersTo = "=Items!`$A`$1:`$G`$19"
    }
}

# ---- Column D width ----------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 12.15

# ---- Selection ------------------------------------------------------
$null = $ws.Range("A9").Select()
